$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update wtkappa (column L) for rows 2-4
$ws.Range("L2").Value = 0.7797696841910529
$ws.Range("L3").Value = 0.7797696841910529
$ws.Range("L4").Value = 0.7821229050279329

# Update SMD (column Q) for rows 2-4
$ws.Range("Q2").Value = -0.02351246133036713
$ws.Range("Q3").Value = -0.02351246133036713
$ws.Range("Q4").Value = 0.01081937260331701
